# Add a new "2022" column (S) to the right of the existing "2021" column (R),
# mirroring the formatting of column R, then move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column R's formatting (number formats, styles, borders, etc.) into
# column S for the header/data rows (4 through 14) so the new column matches
# its neighbour, then fill in the 2022 values.
$ws.Range("R4:R14").Copy()
$ws.Range("S4:S14").PasteSpecial(-4122)  # xlPasteFormats

$values = @{
    4  = 2022
    5  = 99.5
    6  = 99.358544044156048
    7  = 99.400057479522914
    8  = 99.513194978221875
    9  = 99.232429839290006
    10 = 99.453093666824671
    11 = 99.686258104998956
    12 = 99.42525365081228
    13 = 99.561275226674468
    14 = 99.831561216970215
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 19).Value = $values[$row]
}

# Move the active cell / selection to U6 (as recorded in the saved sheet view).
$ws.Range("U6").Select()
